# Auto-applies the cryptos.xlsx crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.831.67'
$ws.Range("E2").Value = '  +4.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.342.92'
$ws.Range("E3").Value = '  +3.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.35'
$ws.Range("E5").Value = '  +2.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.18'
$ws.Range("E6").Value = '  +4.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.341.72'
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +7.39%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").Value = '  +6.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.68'
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.756.32'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.788.41'
$ws.Range("E16").Value = '  +4.06%  '
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.355.08'
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.42'
$ws.Range("E19").Value = '  +1.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.33'
$ws.Range("E21").Value = '  +4.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.52'
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.64'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.161'
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.76'
$ws.Range("E27").Value = '  +4.26%  '
$ws.Range("E28").Value = '  +10.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.17'
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0741'
$ws.Range("E30").Value = '  +5.97%  '
$ws.Range("E31").Value = '  +4.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.23'
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.24'
$ws.Range("E36").Value = '  +3.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.927'
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.98'
$ws.Range("E38").Value = '  +5.51%  '
$ws.Range("E39").Value = '  +7.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.70'
$ws.Range("E40").Value = '  +3.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.379'
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  +5.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.60'
$ws.Range("E43").Value = '  +7.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '275.59'
$ws.Range("E44").Value = '  +10.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.06'
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0932'
$ws.Range("E46").Value = '  +3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0503'
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  +4.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.379'
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.80'
$ws.Range("E51").Value = '  +2.48%  '
